$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "58.80" -> 58.8, or "20.00" -> 20) are first forced to Text format
# so the literal string from the source data is preserved exactly.

$ws.Range("D2").Value = "42.443.87"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.184.63"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.19"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.00"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").Value = "  -5.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.28"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.78"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").Value = "2.510.53"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.10"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").Value = "2.191.06"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.769"
$ws.Range("E17").Value = "  -5.52%  "
$ws.Range("D18").Value = "42.366.72"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  -3.62%  "
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.71"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.41"
$ws.Range("E23").Value = "  -10.09%  "
$ws.Range("E24").Value = "  -5.30%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.45"
$ws.Range("E26").Value = "  -4.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.38"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.03"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.16"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.64"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.02"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0821"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0339"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.22"
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.94"
$ws.Range("E39").Value = "  -9.65%  "
$ws.Range("E40").Value = "  -3.89%  "
$ws.Range("E41").Value = "  +11.90%  "
$ws.Range("E42").Value = "  -7.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.80"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.26"
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.461"
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("E48").Value = "  -4.97%  "
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("E50").Value = "  -2.59%  "
